$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "`$`{`{ deposit_amount `}`}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{{ currency(deposit_amount) }}",
    2
)
